$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2: empty -> "June"
$ws.Range("F2").Value = "June"

# F3: empty -> "'084S"
# A leading apostrophe typed directly into a cell is treated by Excel as a
# "quote prefix" (text-format marker) and gets stripped from the stored
# value. To store the literal apostrophe character as real text content we
# compute it via a formula in a scratch cell and paste the *value* into the
# target cell, which bypasses the quote-prefix interpretation.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = "=""'084S"""
$scratch.Copy()
$ws.Range("F3").PasteSpecial(-4163)
$scratch.ClearContents()

# F4: empty -> "]"
$ws.Range("F4").Value = "]"

# F5: "11/02/2000" -> empty
$ws.Range("F5").Value = ""

$excel.CutCopyMode = $false
